# Apply "[Fonds de solidarite] Add 2022-06-24 data" update:
# refresh nombre_aides (col C) and montant_total (col E) for the affected NAF-section rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 278215
$ws.Range("E10").Value = 1752682990

$ws.Range("C107").Value = 26898
$ws.Range("E107").Value = 36784412

$ws.Range("C167").Value = 101542
$ws.Range("E167").Value = 195302952

$ws.Range("C168").Value = 285128
$ws.Range("E168").Value = 1213925022

$ws.Range("C169").Value = 562684
$ws.Range("E169").Value = 1286541834

$ws.Range("C170").Value = 367606
$ws.Range("E170").Value = 2848481747

$ws.Range("C171").Value = 115231
$ws.Range("E171").Value = 449035067

$ws.Range("C173").Value = 54398
$ws.Range("E173").Value = 151968108

$ws.Range("C174").Value = 357395
$ws.Range("E174").Value = 1020369106

$ws.Range("C175").Value = 125779
$ws.Range("E175").Value = 816694666

$ws.Range("C177").Value = 96785
$ws.Range("E177").Value = 174820350

$ws.Range("C179").Value = 235818
$ws.Range("E179").Value = 813768825

$ws.Range("C180").Value = 141534
$ws.Range("E180").Value = 341253628

$ws.Range("C210").Value = 6433
$ws.Range("E210").Value = 19996314

$ws.Range("C266").Value = 71668
$ws.Range("E266").Value = 219463224

$ws.Range("C279").Value = 28969
$ws.Range("E279").Value = 57090504

$ws.Range("C312").Value = 75105
$ws.Range("E312").Value = 201408498

$ws.Range("C313").Value = 220665
$ws.Range("E313").Value = 1371216867

$ws.Range("C317").Value = 103596
$ws.Range("E317").Value = 303435659

$wb.Save()
